$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")

# ALC row 41
$wsALC.Range("H41").Value = 478.2
$wsALC.Range("J41").Value = 545
$wsALC.Range("L41").Value = 545
$wsALC.Range("N41").Value = -1425

# ALC row 53
$wsALC.Range("H53").Value = 8717.571
$wsALC.Range("I53").Value = 13350.444
$wsALC.Range("K53").Value = 13350.444
$wsALC.Range("M53").Value = -12713.444

# ALC row 64
$wsALC.Range("H64").Value = 3000
$wsALC.Range("I64").Value = 3000
$wsALC.Range("K64").Value = 3000
$wsALC.Range("M64").Value = -2752

# ALC row 67
$wsALC.Range("H67").Value = 3000
$wsALC.Range("I67").Value = 3000
$wsALC.Range("K67").Value = 3000
$wsALC.Range("M67").Value = -2142

# ALC row 74
$wsALC.Range("H74").Value = 3779
$wsALC.Range("I74").Value = 3499
$wsALC.Range("K74").Value = 3499
$wsALC.Range("M74").Value = -2563

# ALC row 77
$wsALC.Range("H77").Value = 3779
$wsALC.Range("I77").Value = 3499
$wsALC.Range("K77").Value = 17495
$wsALC.Range("M77").Value = -12815

# ALC row 92
$wsALC.Range("H92").Value = 2052169.9
$wsALC.Range("I92").Value = 2052169.9
$wsALC.Range("J92").Value = 0
$wsALC.Range("K92").Value = 2052169.9
$wsALC.Range("L92").Value = 0
$wsALC.Range("M92").Value = -2050921.9
$wsALC.Range("N92").ClearContents()

# ALC row 125
$wsALC.Range("H125").Value = 1045.0834
$wsALC.Range("I125").Value = 1062.8
$wsALC.Range("K125").Value = 9565.199999999999
$wsALC.Range("M125").Value = -7105.199999999999

# ALC row 132
$wsALC.Range("H132").Value = 1071.3429
$wsALC.Range("I132").Value = 1081.7273
$wsALC.Range("K132").Value = 3245.1819
$wsALC.Range("M132").Value = -715.1819

# ALC row 137
$wsALC.Range("H137").Value = 1793.95
$wsALC.Range("I137").Value = 1386
$wsALC.Range("J137").Value = 2013.6154
$wsALC.Range("K137").Value = 4158
$wsALC.Range("L137").Value = 6040.8462
$wsALC.Range("M137").Value = -1608
$wsALC.Range("N137").Value = -11140.8462

# ALC row 139
$wsALC.Range("H139").Value = 47583.332
$wsALC.Range("J139").Value = 47583.332
$wsALC.Range("L139").Value = 47583.332
$wsALC.Range("N139").Value = -57863.332

# ARM row 13
$wsARM.Range("H13").Value = 7400000
$wsARM.Range("I13").Value = 7400000
$wsARM.Range("K13").Value = 7400000
$wsARM.Range("M13").Value = -7399856

# ARM row 32
$wsARM.Range("H32").Value = 3867.6296
$wsARM.Range("I32").Value = 2492.1667
$wsARM.Range("K32").Value = 2492.1667
$wsARM.Range("M32").Value = -2205.1667

# ARM row 74
$wsARM.Range("H74").Value = 499.4
$wsARM.Range("I74").Value = 499.4
$wsARM.Range("K74").Value = 499.4
$wsARM.Range("M74").Value = 374.6

# ARM row 77
$wsARM.Range("H77").Value = 499.4
$wsARM.Range("I77").Value = 499.4
$wsARM.Range("K77").Value = 2497
$wsARM.Range("M77").Value = 1871

# ARM row 102
$wsARM.Range("H102").Value = 0
$wsARM.Range("I102").Value = 0
$wsARM.Range("K102").Value = 0
$wsARM.Range("M102").ClearContents()

# BSM row 134
$wsBSM.Range("H134").Value = 7837.5947
$wsBSM.Range("I134").Value = 9139.538
$wsBSM.Range("K134").Value = 27418.614
$wsBSM.Range("M134").Value = -24883.614

# CRP row 7
$wsCRP.Range("H7").Value = 99.14286
$wsCRP.Range("I7").Value = 124.8
$wsCRP.Range("J7").Value = 35
$wsCRP.Range("K7").Value = 124.8
$wsCRP.Range("L7").Value = 35
$wsCRP.Range("M7").Value = -11.8
$wsCRP.Range("N7").Value = -261

# CRP row 31
$wsCRP.Range("H31").Value = 4186.4287
$wsCRP.Range("I31").Value = 1715.7142
$wsCRP.Range("J31").Value = 6657.143
$wsCRP.Range("K31").Value = 1715.7142
$wsCRP.Range("L31").Value = 6657.143
$wsCRP.Range("M31").Value = -1420.7142
$wsCRP.Range("N31").Value = -7247.143

# CRP row 34
$wsCRP.Range("H34").Value = 4186.4287
$wsCRP.Range("I34").Value = 1715.7142
$wsCRP.Range("J34").Value = 6657.143
$wsCRP.Range("K34").Value = 1715.7142
$wsCRP.Range("L34").Value = 6657.143
$wsCRP.Range("M34").Value = -1513.7142
$wsCRP.Range("N34").Value = -7061.143

# CRP row 132
$wsCRP.Range("H132").Value = 2207.3333
$wsCRP.Range("I132").Value = 1146.4445
$wsCRP.Range("K132").Value = 3439.3335
$wsCRP.Range("M132").Value = -909.3335000000002

# CRP row 134
$wsCRP.Range("H134").Value = 774.1818
$wsCRP.Range("I134").Value = 651.6
$wsCRP.Range("K134").Value = 1954.8
$wsCRP.Range("M134").Value = 580.1999999999998

# CUL row 92
$wsCUL.Range("H92").Value = 325
$wsCUL.Range("J92").Value = 325
$wsCUL.Range("L92").Value = 975
$wsCUL.Range("N92").Value = -3471

# CUL row 113
$wsCUL.Range("H113").Value = 7493.6665
$wsCUL.Range("I113").Value = 25576.75
$wsCUL.Range("J113").Value = 918
$wsCUL.Range("K113").Value = 76730.25
$wsCUL.Range("L113").Value = 2754
$wsCUL.Range("M113").Value = -74560.25
$wsCUL.Range("N113").Value = -7094

# CUL row 129
$wsCUL.Range("H129").Value = 64474.363
$wsCUL.Range("J129").Value = 117700.336
$wsCUL.Range("L129").Value = 353101.008
$wsCUL.Range("N129").Value = -363101.008

# GSM row 80
$wsGSM.Range("H80").Value = 2598.75
$wsGSM.Range("I80").Value = 2435
$wsGSM.Range("J80").Value = 2653.3333
$wsGSM.Range("K80").Value = 2435
$wsGSM.Range("L80").Value = 2653.3333
$wsGSM.Range("M80").Value = -1437
$wsGSM.Range("N80").Value = -4649.3333

# GSM row 83
$wsGSM.Range("H83").Value = 2598.75
$wsGSM.Range("I83").Value = 2435
$wsGSM.Range("J83").Value = 2653.3333
$wsGSM.Range("K83").Value = 12175
$wsGSM.Range("L83").Value = 13266.6665
$wsGSM.Range("M83").Value = -7183
$wsGSM.Range("N83").Value = -23250.6665

# GSM row 97
$wsGSM.Range("H97").Value = 0
$wsGSM.Range("I97").Value = 0
$wsGSM.Range("J97").Value = 0
$wsGSM.Range("K97").Value = 0
$wsGSM.Range("L97").Value = 0
$wsGSM.Range("M97").ClearContents()
$wsGSM.Range("N97").ClearContents()

# GSM row 102
$wsGSM.Range("H102").Value = 2154.724
$wsGSM.Range("I102").Value = 2127.389
$wsGSM.Range("K102").Value = 2127.389
$wsGSM.Range("M102").Value = -505.3890000000001

# GSM row 132
$wsGSM.Range("H132").Value = 5289.0454
$wsGSM.Range("I132").Value = 4265.647
$wsGSM.Range("K132").Value = 12796.941
$wsGSM.Range("M132").Value = -10266.941

# LTW row 7
$wsLTW.Range("H7").Value = 8612.625
$wsLTW.Range("J7").Value = 8413.857
$wsLTW.Range("L7").Value = 8413.857
$wsLTW.Range("N7").Value = -8637.857

# LTW row 12
$wsLTW.Range("H12").Value = 0
$wsLTW.Range("I12").Value = 0
$wsLTW.Range("K12").Value = 0
$wsLTW.Range("M12").ClearContents()

# LTW row 82
$wsLTW.Range("H82").Value = 4868.4287
$wsLTW.Range("I82").Value = 2000
$wsLTW.Range("K82").Value = 2000
$wsLTW.Range("M82").Value = -1639

# LTW row 85
$wsLTW.Range("H85").Value = 4868.4287
$wsLTW.Range("I85").Value = 2000
$wsLTW.Range("K85").Value = 2000
$wsLTW.Range("M85").Value = -752

# LTW row 126
$wsLTW.Range("H126").Value = 8612.625
$wsLTW.Range("J126").Value = 8413.857
$wsLTW.Range("L126").Value = 25241.571
$wsLTW.Range("N126").Value = -30181.571

# LTW row 132
$wsLTW.Range("H132").Value = 1988.091
$wsLTW.Range("I132").Value = 1120
$wsLTW.Range("K132").Value = 3360
$wsLTW.Range("M132").Value = -830

# LTW row 136
$wsLTW.Range("H136").Value = 3214.0454
$wsLTW.Range("I136").Value = 1929.4
$wsLTW.Range("K136").Value = 5788.200000000001
$wsLTW.Range("M136").Value = -3238.200000000001
